$d = $word.ActiveDocument

for ($i = 1; $i -le 6; $i++) {
    $old = "<id>p053r_a$i</id>"
    $new = "<id>p053r_$i</id>"
    $found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced #$i :" $found
}
